# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets.
# Each row has: Date | Timestamp | Hour | Location | Value | Status

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param($SheetName, $StartRow, $Data)

    $ws = $wb.Worksheets.Item($SheetName)

    $lines = ($Data.Trim() -split "`n")
    $nRows = $lines.Count
    $endRow = $StartRow + $nRows - 1

    # Plain ".Value = ..." lets Excel auto-detect dates/percentages/numbers
    # (e.g. "2026-02-01" -> a real date, "77.9%" -> 0.779 with a percent
    # style) which would add number-format styles not present in the
    # target file. Writing every cell as a `="literal"` formula first
    # keeps it as plain text, then Copy + PasteSpecial (values only) bakes
    # the whole block back down to plain literal string cells with no
    # leftover formulas and no leftover styles.
    for ($i = 0; $i -lt $nRows; $i++) {
        $parts = ($lines[$i]).Trim() -split '\|'
        $r = $StartRow + $i

        $ws.Range("A$r").Formula = "=""2026-02-01"""
        $ws.Range("B$r").Formula = "=""" + $parts[0] + """"
        $ws.Range("C$r").Formula = "=""14:00"""
        $ws.Range("D$r").Formula = "=""Bathroom"""
        $ws.Range("E$r").Formula = "=""" + $parts[1] + """"
        $ws.Range("F$r").Formula = "=""" + $parts[2] + """"
    }

    $rng = $ws.Range("A$StartRow" + ":F$endRow")
    $rng.Copy()
    $rng.PasteSpecial(-4163)
    $excel.CutCopyMode = 0
}

# PIR sheet: append rows 241-254 (Timestamp|Value|Status)
$pirData = @"
14:13:50|No Motion|Inactive
14:13:52|No Motion|Inactive
14:13:53|No Motion|Inactive
14:13:58|No Motion|Inactive
14:14:03|No Motion|Inactive
14:14:08|No Motion|Inactive
14:14:13|No Motion|Inactive
14:14:18|No Motion|Inactive
14:14:23|No Motion|Inactive
14:14:28|No Motion|Inactive
14:14:33|No Motion|Inactive
14:14:38|No Motion|Inactive
14:14:43|No Motion|Inactive
14:14:48|No Motion|Inactive
"@
Add-LogRows "PIR" 241 $pirData

# Humidity sheet: append rows 158-170 (Timestamp|Value|Status)
$humidityData = @"
14:13:49|77.9%|Active
14:13:51|76.9%|Active
14:13:52|77.9%|Active
14:13:56|77.0%|Active
14:14:01|77.9%|Active
14:14:06|77.0%|Active
14:14:11|78.0%|Active
14:14:17|76.9%|Active
14:14:22|78.0%|Active
14:14:26|76.9%|Active
14:14:32|77.9%|Active
14:14:42|77.9%|Active
14:14:47|76.9%|Active
"@
Add-LogRows "Humidity" 158 $humidityData

# Temperature sheet: append rows 79-91 (Timestamp|Value|Status)
$temperatureData = @"
14:13:50|29.4C|Active
14:13:51|29.4C|Active
14:13:53|29.4C|Active
14:13:57|29.5C|Active
14:14:02|29.4C|Active
14:14:07|29.5C|Active
14:14:12|29.5C|Active
14:14:17|29.4C|Active
14:14:22|29.4C|Active
14:14:27|29.4C|Active
14:14:32|29.5C|Active
14:14:42|29.4C|Active
14:14:47|29.4C|Active
"@
Add-LogRows "Temperature" 79 $temperatureData
